# Apply the odds updates described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "N3"  = 5
    "O3"  = 1.67
    "P3"  = 2.1
    "S3"  = 1.73
    "T3"  = 2.08
    "U3"  = 2.75
    "V3"  = 1.4
    "AA3" = 21
    "AC3" = 5
    "AZ3" = 301

    "G4"  = 3.05
    "H4"  = 3.05
    "I4"  = 2.3
    "J4"  = 3.5
    "K4"  = 2.07
    "L4"  = 2.85
    "N4"  = 8.1
    "AA4" = 26
    "AG4" = 7.7
    "AK4" = 19
    "AL4" = 28
    "AO4" = 16
    "AP4" = 21
    "AQ4" = 75
    "AS4" = 2.6
    "AV4" = 4.25
    "AW4" = 11.75
    "AX4" = 18
    "AY4" = 45
    "BA4" = 250
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
